# Update of all scripts and data
# The source dataset dropped the "Mullus barbatus" catch record (row 9) for
# station 75 / gear 1-RAP, so every subsequent record shifts up by one row,
# and the relative-frequency (RF / column I) figures for the 2-RAP discard
# block are recalculated against the new totals. Also, items that used to be
# tallied with a count (Numb / column H) of 0 are now tallied as -1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Mullus barbatus" row entirely - this shifts every row below it
# up by one and keeps all the other column data (Survey/Area/Station/Gear)
# intact automatically.
$ws.Rows.Item(9).Delete()

# After the shift, fix up the count (column H) for the "no specimen counted"
# discard-type rows, which change from 0 to -1.
$zeroCountRows = @(33, 36, 40, 45)
foreach ($r in $zeroCountRows) {
    $ws.Cells.Item($r, 8).Value = -1
}

# Update the relative frequency (column I) values for the 2-RAP discard
# block, which were recalculated for the new dataset.
$ws.Cells.Item(30, 9).Value = 5.68325
for ($r = 31; $r -le 45; $r++) {
    $ws.Cells.Item($r, 9).Value = 10.3665
}
